# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback DateTime" (H2)
# timestamps on the zh-cn and de-de sheets to reflect a re-run of the handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-13 11:05:27"
$wsZhCn.Range("H2").Value = "2016-03-13 11:05:44"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-13 11:05:31"
$wsDeDe.Range("H2").Value = "2016-03-13 11:05:50"
